# Apply the data-rotation edit to rows 24-27 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 24 -> becomes the old "Rosenticka" record (with new location text) ----
$ws.Range("A24").Value = 111957798
$ws.Range("B24").Value = 89820
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 658
$ws.Range("F24").Value = "Rosenticka"
$ws.Range("G24").Value = "Rhodofomes roseus"
$ws.Range("H24").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("I24").Value = "'6"
$ws.Range("J24").Value = "fruktkroppar"
$ws.Range("P24").Value = "Österåsen, Österås, Ång"
$ws.Range("Q24").Value = 609747
$ws.Range("R24").Value = 7011953
# L24 and AC24 cells are removed in the target workbook
$ws.Range("L24").ClearContents()
$ws.Range("AC24").ClearContents()

# ---- Row 25 -> keeps the same species data, only id/location shift ----
$ws.Range("A25").Value = 111957843
$ws.Range("B25").Value = 89820
$ws.Range("P25").Value = "Österås, Österås, Ång"
$ws.Range("Q25").Value = 609773
$ws.Range("R25").Value = 7011992

# ---- Row 26 -> becomes the old "Knärot" record ----
$ws.Range("A26").Value = 111958205
$ws.Range("B26").Value = 96720
$ws.Range("D26").Value = "VU"
$ws.Range("E26").Value = 220787
$ws.Range("F26").Value = "Knärot"
$ws.Range("G26").Value = "Goodyera repens"
$ws.Range("H26").Value = "(L.) R. Br."
$ws.Range("I26").Value = "'1"
$ws.Range("J26").Value = "plantor/tuvor"
$ws.Range("P26").Value = "Österåsen, Ång"
$ws.Range("Q26").Value = 609803
$ws.Range("R26").Value = 7011969
$ws.Range("AC26").Value = "½ m2"

# ---- Row 27 -> only the taxon sort order changes ----
$ws.Range("B27").Value = 55643
